# "moved celsius to be the default" — add three new Air Temperature fields
# (Max/Min/WeightedMean) to the gc_fields_uom sheet, with the metric (C)
# column set to celsius and the statute (D) column set to fahrenheit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gc_fields_uom")

# Populate the new rows in an order that reproduces the original shared-
# string insertion sequence (field name, then the "fahrenheit" unit, then
# the remaining field names, then "celsius", then the repeated values).
$ws.Range("A64").Value = "MaxAirTemperature"
$ws.Range("D64").Value = "fahrenheit"
$ws.Range("A65").Value = "MinAirTemperature"
$ws.Range("A66").Value = "WeightedMeanAirTemperature"
$ws.Range("C64").Value = "celsius"

$ws.Range("B64").Value = "all"
$ws.Range("B65").Value = "all"
$ws.Range("B66").Value = "all"

$ws.Range("C65").Value = "celsius"
$ws.Range("C66").Value = "celsius"

$ws.Range("D65").Value = "fahrenheit"
$ws.Range("D66").Value = "fahrenheit"

# Give the new rows the (new) 11pt default font style used for the rest of
# the added data.
$ws.Range("A64:D66").Font.Size = 11

# Leave the selection on the last entered cell, matching the saved view.
$ws.Range("C66").Select()
